$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 2).Value = 0.2028864896758006
$ws.Cells.Item(1, 3).Value = 0.002539495620456482
$ws.Cells.Item(1, 4).Value = -1.379316785363834
$ws.Cells.Item(1, 5).Value = 0.188939888900502
$ws.Cells.Item(1, 6).Value = 1.570796384046431
$ws.Cells.Item(2, 2).Value = 0.2046333296544263
$ws.Cells.Item(2, 3).Value = 0.002533299447576215
$ws.Cells.Item(2, 4).Value = -1.379362809194377
$ws.Cells.Item(2, 5).Value = 0.1889000613375225
$ws.Cells.Item(2, 6).Value = 1.570796384060886
$ws.Cells.Item(3, 2).Value = 0.215527422831531
$ws.Cells.Item(3, 3).Value = 0.002494657279132034
$ws.Cells.Item(3, 4).Value = -1.379649834854683
$ws.Cells.Item(3, 5).Value = 0.1886516784361557
$ws.Cells.Item(3, 6).Value = 1.570796384151034
$ws.Cells.Item(4, 2).Value = 0.2412914350012368
$ws.Cells.Item(4, 3).Value = 0.002403270388835187
$ws.Cells.Item(4, 4).Value = -1.380328636888492
$ws.Cells.Item(4, 5).Value = 0.1880642646891343
$ws.Cells.Item(4, 6).Value = 1.570796384364231
$ws.Cells.Item(5, 2).Value = 0.2844036859956439
$ws.Cells.Item(5, 3).Value = 0.00225034798971023
$ws.Cells.Item(5, 4).Value = -1.381464511358056
$ws.Cells.Item(5, 5).Value = 0.1870813149555161
$ws.Cells.Item(5, 6).Value = 1.570796384720984
$ws.Cells.Item(6, 2).Value = 0.3446388740118342
$ws.Cells.Item(6, 3).Value = 0.002036689244209476
$ws.Cells.Item(6, 4).Value = -1.383051522257715
$ws.Cells.Item(6, 5).Value = 0.1857079660662958
$ws.Cells.Item(6, 6).Value = 1.570796385219429
$ws.Cells.Item(7, 2).Value = 0.4196087999388747
$ws.Cells.Item(7, 3).Value = 0.001770765274327443
$ws.Cells.Item(7, 4).Value = -1.385026747927482
$ws.Cells.Item(7, 5).Value = 0.1839986684300178
$ws.Cells.Item(7, 6).Value = 1.570796385839804
$ws.Cells.Item(8, 2).Value = 0.5053030916848213
$ws.Cells.Item(8, 3).Value = 0.001466801171715305
$ws.Cells.Item(8, 4).Value = -1.387284527466625
$ws.Cells.Item(8, 5).Value = 0.1820448576383891
$ws.Cells.Item(8, 6).Value = 1.570796386548922
$ws.Cells.Item(9, 2).Value = 0.5966299285037232
$ws.Cells.Item(9, 3).Value = 0.001142858007795339
$ws.Cells.Item(9, 4).Value = -1.389690707147247
$ws.Cells.Item(9, 5).Value = 0.1799626260718912
$ws.Cells.Item(9, 6).Value = 1.57079638730465
$ws.Cells.Item(10, 2).Value = 0.6879567653226248
$ws.Cells.Item(10, 3).Value = 0.0008189148438753735
$ws.Cells.Item(10, 4).Value = -1.392096886827868
$ws.Cells.Item(10, 5).Value = 0.1778803945053933
$ws.Cells.Item(10, 6).Value = 1.570796388060378
$ws.Cells.Item(11, 2).Value = 0.7736510570685715
$ws.Cells.Item(11, 3).Value = 0.0005149507412632363
$ws.Cells.Item(11, 4).Value = -1.394354666367011
$ws.Cells.Item(11, 5).Value = 0.1759265837137646
$ws.Cells.Item(11, 6).Value = 1.570796388769496
$ws.Cells.Item(12, 2).Value = 0.8486209829956121
$ws.Cells.Item(12, 3).Value = 0.0002490267713812027
$ws.Cells.Item(12, 4).Value = -1.396329892036778
$ws.Cells.Item(12, 5).Value = 0.1742172860774866
$ws.Cells.Item(12, 6).Value = 1.570796389389871
$ws.Cells.Item(13, 2).Value = 0.908856171011802
$ws.Cells.Item(13, 3).Value = 0.0000353680258804488570346
$ws.Cells.Item(13, 4).Value = -1.397916902936438
$ws.Cells.Item(13, 5).Value = 0.1728439371882663
$ws.Cells.Item(13, 6).Value = 1.570796389888316
$ws.Cells.Item(14, 2).Value = 0.9519684220062091
$ws.Cells.Item(14, 3).Value = -0.0001175543732445084
$ws.Cells.Item(14, 4).Value = -1.399052777406001
$ws.Cells.Item(14, 5).Value = 0.1718609874546481
$ws.Cells.Item(14, 6).Value = 1.570796390245069
$ws.Cells.Item(15, 2).Value = 0.9777324341759148
$ws.Cells.Item(15, 3).Value = -0.0002089412635413557
$ws.Cells.Item(15, 4).Value = -1.39973157943981
$ws.Cells.Item(15, 5).Value = 0.1712735737076268
$ws.Cells.Item(15, 6).Value = 1.570796390458266
$ws.Cells.Item(16, 2).Value = 0.9886265273530191
$ws.Cells.Item(16, 3).Value = -0.0002475834319855381
$ws.Cells.Item(16, 4).Value = -1.400018605100116
$ws.Cells.Item(16, 5).Value = 0.1710251908062599
$ws.Cells.Item(16, 6).Value = 1.570796390548414
$ws.Cells.Item(17, 2).Value = 0.9903733673316449
$ws.Cells.Item(17, 3).Value = -0.0002537796048658041
$ws.Cells.Item(17, 4).Value = -1.40006462893066
$ws.Cells.Item(17, 5).Value = 0.1709853632432805
$ws.Cells.Item(17, 6).Value = 1.570796390562869

Write-Output "Applied robot trajectory sign/angle fixes"
